$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '22.057.90'
$ws.Range("E2").Value = '  -1.13%  '

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '1.556.39'
$ws.Range("E3").Value = '  -0.32%  '

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '1.002'
$ws.Range("E4").Value = '  +0.16%  '

$ws.Range("E5").Value = '  +0.07%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '291.61'
$ws.Range("E6").Value = '  +0.62%  '

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.3927'
$ws.Range("E7").Value = '  +4.09%  '

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.3230'

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '43.75'
$ws.Range("E9").Value = '  -2.20%  '

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.07301'
$ws.Range("E10").Value = '  -1.34%  '

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '1.084'
$ws.Range("E11").Value = '  -5.54%  '

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '1.002'
$ws.Range("E12").Value = '  +0.14%  '

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '18.87'
$ws.Range("E13").Value = '  -7.37%  '

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '5.685'
$ws.Range("E14").Value = '  -3.43%  '

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '0.00001128'
$ws.Range("E15").Value = '  +4.32%  '

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '6.640'
$ws.Range("E16").Value = '  -2.07%  '

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '1.554.22'
$ws.Range("E17").Value = '  -0.94%  '

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '0.06588'
$ws.Range("E18").Value = '  -1.08%  '

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '83.77'
$ws.Range("E19").Value = '  -3.19%  '

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '1.001'
$ws.Range("E20").Value = '  +0.06%  '

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '6.302'
$ws.Range("E21").Value = '  -2.15%  '

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '15.69'
$ws.Range("E22").Value = '  -3.34%  '

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '11.32'
$ws.Range("E23").Value = '  -3.76%  '

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '22.071.58'
$ws.Range("E24").Value = '  -1.04%  '

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '2.355'
$ws.Range("E25").Value = '  +2.49%  '

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '2.445'
$ws.Range("E26").Value = '  -5.95%  '

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '149.02'
$ws.Range("E27").Value = '  -1.47%  '

$ws.Range("E28").Value = '  -3.73%  '

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '4.893'
$ws.Range("E29").Value = '  -0.85%  '

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '1.728.10'
$ws.Range("E30").Value = '  -0.84%  '

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '119.29'
$ws.Range("E31").Value = '  -3.32%  '

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '1.009'
$ws.Range("E32").Value = '  -7.32%  '

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '5.808'
$ws.Range("E33").Value = '  -2.54%  '

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '0.08326'
$ws.Range("E34").Value = '  +1.44%  '

$ws.Range("E35").Value = '  -15.93%  '

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '9.094'
$ws.Range("E36").Value = '  -4.18%  '

$ws.Range("B37").Value = 'VeChain'
$ws.Range("C37").Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '0.02277'
$ws.Range("E37").Value = '  -3.88%  '

$ws.Range("B38").Value = 'Hedera'
$ws.Range("C38").Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.06116'
$ws.Range("E38").Value = '  -4.00%  '

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '5.145'
$ws.Range("E39").Value = '  -4.58%  '

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '1.213'
$ws.Range("E40").Value = '  -3.09%  '

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.2055'
$ws.Range("E41").Value = '  -5.43%  '

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '1.001'
$ws.Range("E42").Value = '  +0.09%  '

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '10.78'
$ws.Range("E43").Value = '  -2.74%  '

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.5847'
$ws.Range("E44").Value = '  -4.26%  '

$ws.Range("B45").Value = 'EnergySwap'
$ws.Range("C45").Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '13.12'
$ws.Range("E45").Value = '  -5.26%  '

$ws.Range("B46").Value = 'PancakeSwap'
$ws.Range("C46").Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '3.759'
$ws.Range("E46").Value = '  -0.19%  '

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '0.5606'
$ws.Range("E47").Value = '  -5.53%  '

$ws.Range("B48").Value = 'Quant'
$ws.Range("C48").Value = 'https://coinranking.com/coin/bauj_21eYVwso+quant-qnt'
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '118.93'
$ws.Range("E48").Value = '  -3.86%  '

$ws.Range("B49").Value = 'NEARProtocol'
$ws.Range("C49").Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '1.909'
$ws.Range("E49").Value = '  -4.18%  '

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '1.138'
$ws.Range("E50").Value = '  -3.77%  '

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.06833'
$ws.Range("E51").Value = '  -3.82%  '
